$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so that numeric-looking
# price strings (e.g. "234.17", "0.600") are preserved verbatim as
# text, matching the inlineStr cells in the original workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "36.564.71"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "2.000.84"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "234.17"
$ws.Range("E5").Value = "  -9.15%  "

$ws.Range("D6").Value = "0.600"
$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "55.20"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").Value = "0.372"
$ws.Range("E9").Value = "  -3.38%  "

$ws.Range("D10").Value = "58.10"
$ws.Range("E10").Value = "  +2.87%  "

$ws.Range("D12").Value = "0.0990"
$ws.Range("E12").Value = "  -3.08%  "

$ws.Range("D13").Value = "14.27"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "2.292.65"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "20.46"
$ws.Range("E15").Value = "  -3.18%  "

$ws.Range("D16").Value = "0.760"
$ws.Range("E16").Value = "  -5.02%  "

$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").Value = "2.001.01"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").Value = "36.540.40"
$ws.Range("E19").Value = "  -1.61%  "

$ws.Range("D20").Value = "67.82"
$ws.Range("E20").Value = "  -2.46%  "

$ws.Range("D21").Value = "0.0₃0806"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").Value = "5.33"
$ws.Range("E22").Value = "  +4.03%  "

$ws.Range("D23").Value = "222.43"
$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -7.96%  "

$ws.Range("D27").Value = "162.16"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").Value = "8.66"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").Value = "0.128"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "18.85"
$ws.Range("E30").Value = "  -3.89%  "

$ws.Range("D31").Value = "1.33"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").Value = "4.39"
$ws.Range("E33").Value = "  -5.01%  "

$ws.Range("D34").Value = "0.0605"
$ws.Range("E34").Value = "  -5.68%  "

$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  -6.01%  "

$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  -3.24%  "

$ws.Range("D40").Value = "5.66"
$ws.Range("E40").Value = "  +7.79%  "

$ws.Range("E41").Value = "  -1.14%  "

$ws.Range("E42").Value = "  +1.79%  "

$ws.Range("D43").Value = "1.457.58"
$ws.Range("E43").Value = "  +4.17%  "

$ws.Range("E44").Value = "  -3.84%  "

$ws.Range("E45").Value = "  -7.83%  "

$ws.Range("D46").Value = "89.37"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "15.24"
$ws.Range("E47").Value = "  -2.94%  "

$ws.Range("D48").Value = "0.996"
$ws.Range("E48").Value = "  -2.39%  "

$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("D50").Value = "6.85"
$ws.Range("E50").Value = "  -2.52%  "

$ws.Range("D51").Value = "3.74"
$ws.Range("E51").Value = "  +8.88%  "

# Restore the default (Normal) style so the cells do not carry an
# explicit/custom style index, matching the original formatting.
$ws.Range("D2:D51").Style = "Normal"
